$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to re-pulled data
$ws.Range("F2").Value = -3
$ws.Range("F3").Value = -3
$ws.Range("F6").Value = -4
$ws.Range("F7").Value = -2
